$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.437.45"
$ws.Range("E2").Value = "  -1.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.844.26"
$ws.Range("E3").Value = "  -1.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "264.72"
$ws.Range("E5").Value = "  -3.47%  "

$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5203"
$ws.Range("E7").Value = "  -1.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3263"
$ws.Range("E8").Value = "  -3.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06798"
$ws.Range("E9").Value = "  -0.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.77"
$ws.Range("E10").Value = "  -5.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7782"
$ws.Range("E11").Value = "  -1.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07745"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.843.92"
$ws.Range("E13").Value = "  -0.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.07"
$ws.Range("E14").Value = "  -2.22%  "

$ws.Range("E15").Value = "  -2.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9993"
$ws.Range("E16").Value = "  +0.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.93"
$ws.Range("E17").Value = "  -3.50%  "

$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007961"
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.448.01"
$ws.Range("E20").Value = "  -1.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.077.29"
$ws.Range("E21").Value = "  -1.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.623"
$ws.Range("E22").Value = "  -1.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.582"
$ws.Range("E23").Value = "  -3.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.996"
$ws.Range("E24").Value = "  -1.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.64"
$ws.Range("E25").Value = "  -1.32%  "

$ws.Range("E26").Value = "  -8.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.651"
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.99"
$ws.Range("E28").Value = "  -1.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.76"
$ws.Range("E29").Value = "  -0.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.165"
$ws.Range("E30").Value = "  -3.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.125"
$ws.Range("E31").Value = "  -4.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08705"
$ws.Range("E32").Value = "  -1.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04826"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.131"
$ws.Range("E34").Value = "  -2.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7199"
$ws.Range("E35").Value = "  -1.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.843"
$ws.Range("E36").Value = "  -1.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.102"
$ws.Range("E37").Value = "  -3.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01781"
$ws.Range("E38").Value = "  -3.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.231"
$ws.Range("E39").Value = "  -4.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4868"
$ws.Range("E40").Value = "  -4.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9148"
$ws.Range("E41").Value = "  -2.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.73"
$ws.Range("E42").Value = "  -4.78%  "

$ws.Range("E43").Value = "  -0.94%  "

$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.729"
$ws.Range("E45").Value = "  -3.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4175"
$ws.Range("E46").Value = "  -5.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05921"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.055"
$ws.Range("E48").Value = "  -2.65%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.08"
$ws.Range("E49").Value = "  -2.81%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1235"
$ws.Range("E50").Value = "  -7.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8856"
$ws.Range("E51").Value = "  +0.61%  "
